# Apply the "AIC_Max" sector-name/ordering fix across every year worksheet.
#
# For every worksheet:
#   1. The values in column E for rows 5, 6, 7 are rotated:
#        new E5 = old E6
#        new E6 = old E7
#        new E7 = old E5
#   2. The Sector labels in column C (rows 5-7) are renamed so the data now
#      aligns with the Baseline naming convention:
#        "Offshore wind" -> "Onshore wind plants"
#        "Onshore wind"  -> "Photovoltaic plants"
#        "PV"            -> "Offshore wind plants"
#
# The renaming is done one column at a time (all sheets for C5, then all
# sheets for C6, then all sheets for C7) so that every sheet's label is
# updated consistently, keeping the data tied to the correct row.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

# --- Step 1: rotate the E5:E7 values on every sheet ------------------------
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $e5 = $ws.Cells.Item(5, 5).Value()
    $e6 = $ws.Cells.Item(6, 5).Value()
    $e7 = $ws.Cells.Item(7, 5).Value()

    $ws.Cells.Item(5, 5).Value = $e6
    $ws.Cells.Item(6, 5).Value = $e7
    $ws.Cells.Item(7, 5).Value = $e5
}

# --- Step 2: rename the Sector labels in column C on every sheet -----------
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Item(5, 3).Value = "Onshore wind plants"
}

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Item(6, 3).Value = "Photovoltaic plants"
}

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Cells.Item(7, 3).Value = "Offshore wind plants"
}
